$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# Use the NumberFormat "@" (Text) trick so numeric-looking strings (e.g. "137.00")
# are stored as text like the original inline strings, then restore the original
# cell style so no stray formatting/style diff is introduced.

$cell = $ws.Range('D2')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '63.119.86'
$cell.Style = $origStyle

$cell = $ws.Range('E2')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -1.87%  '
$cell.Style = $origStyle

$cell = $ws.Range('D3')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '3.134.34'
$cell.Style = $origStyle

$cell = $ws.Range('E3')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.08%  '
$cell.Style = $origStyle

$cell = $ws.Range('E4')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +0.08%  '
$cell.Style = $origStyle

$cell = $ws.Range('D5')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '594.83'
$cell.Style = $origStyle

$cell = $ws.Range('E5')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -2.50%  '
$cell.Style = $origStyle

$cell = $ws.Range('D6')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '137.00'
$cell.Style = $origStyle

$cell = $ws.Range('E6')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -4.39%  '
$cell.Style = $origStyle

$cell = $ws.Range('E7')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +0.09%  '
$cell.Style = $origStyle

$cell = $ws.Range('D8')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '3.120.85'
$cell.Style = $origStyle

$cell = $ws.Range('E8')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.43%  '
$cell.Style = $origStyle

$cell = $ws.Range('D9')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.519'
$cell.Style = $origStyle

$cell = $ws.Range('E9')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -1.86%  '
$cell.Style = $origStyle

$cell = $ws.Range('E10')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -2.66%  '
$cell.Style = $origStyle

$cell = $ws.Range('D11')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '5.29'
$cell.Style = $origStyle

$cell = $ws.Range('E11')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -1.14%  '
$cell.Style = $origStyle

$cell = $ws.Range('D12')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.460'
$cell.Style = $origStyle

$cell = $ws.Range('E12')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -3.07%  '
$cell.Style = $origStyle

$cell = $ws.Range('E13')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -2.51%  '
$cell.Style = $origStyle

$cell = $ws.Range('D14')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '34.36'
$cell.Style = $origStyle

$cell = $ws.Range('E14')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -3.14%  '
$cell.Style = $origStyle

$cell = $ws.Range('D15')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '3.645.62'
$cell.Style = $origStyle

$cell = $ws.Range('E15')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.14%  '
$cell.Style = $origStyle

$cell = $ws.Range('D17')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '63.188.48'
$cell.Style = $origStyle

$cell = $ws.Range('E17')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -1.71%  '
$cell.Style = $origStyle

$cell = $ws.Range('D18')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '3.122.85'
$cell.Style = $origStyle

$cell = $ws.Range('E18')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -1.05%  '
$cell.Style = $origStyle

$cell = $ws.Range('D19')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '6.74'
$cell.Style = $origStyle

$cell = $ws.Range('E19')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -1.59%  '
$cell.Style = $origStyle

$cell = $ws.Range('D20')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '477.90'
$cell.Style = $origStyle

$cell = $ws.Range('E20')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.06%  '
$cell.Style = $origStyle

$cell = $ws.Range('D21')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '14.18'
$cell.Style = $origStyle

$cell = $ws.Range('E21')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -3.45%  '
$cell.Style = $origStyle

$cell = $ws.Range('D22')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.700'
$cell.Style = $origStyle

$cell = $ws.Range('E22')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -2.94%  '
$cell.Style = $origStyle

$cell = $ws.Range('D23')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '7.73'
$cell.Style = $origStyle

$cell = $ws.Range('E23')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.89%  '
$cell.Style = $origStyle

$cell = $ws.Range('D24')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '87.24'
$cell.Style = $origStyle

$cell = $ws.Range('E24')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +2.67%  '
$cell.Style = $origStyle

$cell = $ws.Range('D25')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '13.07'
$cell.Style = $origStyle

$cell = $ws.Range('E25')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -3.88%  '
$cell.Style = $origStyle

$cell = $ws.Range('E26')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +0.11%  '
$cell.Style = $origStyle

$cell = $ws.Range('D27')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.73'
$cell.Style = $origStyle

$cell = $ws.Range('E27')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -1.95%  '
$cell.Style = $origStyle

$cell = $ws.Range('D28')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '7.20'
$cell.Style = $origStyle

$cell = $ws.Range('E29')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -5.98%  '
$cell.Style = $origStyle

$cell = $ws.Range('D30')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.10'
$cell.Style = $origStyle

$cell = $ws.Range('D31')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '27.09'
$cell.Style = $origStyle

$cell = $ws.Range('E31')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +1.54%  '
$cell.Style = $origStyle

$cell = $ws.Range('E32')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +0.05%  '
$cell.Style = $origStyle

$cell = $ws.Range('D33')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.108'
$cell.Style = $origStyle

$cell = $ws.Range('E33')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -7.83%  '
$cell.Style = $origStyle

$cell = $ws.Range('E34')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -3.17%  '
$cell.Style = $origStyle

$cell = $ws.Range('E35')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -2.63%  '
$cell.Style = $origStyle

$cell = $ws.Range('D36')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '5.85'
$cell.Style = $origStyle

$cell = $ws.Range('E36')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -1.49%  '
$cell.Style = $origStyle

$cell = $ws.Range('D37')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '51.92'
$cell.Style = $origStyle

$cell = $ws.Range('E37')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -1.01%  '
$cell.Style = $origStyle

$cell = $ws.Range('D38')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0₃0712'
$cell.Style = $origStyle

$cell = $ws.Range('E38')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -3.79%  '
$cell.Style = $origStyle

$cell = $ws.Range('D39')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0391'
$cell.Style = $origStyle

$cell = $ws.Range('E39')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -1.01%  '
$cell.Style = $origStyle

$cell = $ws.Range('D40')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '420.77'
$cell.Style = $origStyle

$cell = $ws.Range('E40')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -7.15%  '
$cell.Style = $origStyle

$cell = $ws.Range('E41')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.89%  '
$cell.Style = $origStyle

$cell = $ws.Range('E42')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.39%  '
$cell.Style = $origStyle

$cell = $ws.Range('D43')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.69'
$cell.Style = $origStyle

$cell = $ws.Range('E43')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -9.57%  '
$cell.Style = $origStyle

$cell = $ws.Range('D44')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.887.52'
$cell.Style = $origStyle

$cell = $ws.Range('E44')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +1.13%  '
$cell.Style = $origStyle

$cell = $ws.Range('E45')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.07%  '
$cell.Style = $origStyle

$cell = $ws.Range('E46')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -5.36%  '
$cell.Style = $origStyle

$cell = $ws.Range('E47')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.06%  '
$cell.Style = $origStyle

$cell = $ws.Range('D48')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '25.76'
$cell.Style = $origStyle

$cell = $ws.Range('E48')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -2.36%  '
$cell.Style = $origStyle

$cell = $ws.Range('E49')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +0.00%  '
$cell.Style = $origStyle

$cell = $ws.Range('D50')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.28'
$cell.Style = $origStyle

$cell = $ws.Range('E50')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -5.65%  '
$cell.Style = $origStyle

$cell = $ws.Range('D51')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '118.96'
$cell.Style = $origStyle

$cell = $ws.Range('E51')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.95%  '
$cell.Style = $origStyle

